$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write cells in the exact order needed so new shared strings are
# --- interned in the order they appear in the target workbook
# --- (AAABBBCCCDDD, B, C, D, aBc!2@#xYz, a, c, !, @, #, x, z,
# ---  Hello World, e, l, W, o, r, d).

# New "AAABBBCCCDDD" test block (numRows = 4 grid in rows 22-25)
$ws.Range("C7").Value = "AAABBBCCCDDD"
$ws.Range("J23").Value = "B"
$ws.Range("K22").Value = "C"
$ws.Range("M23").Value = "D"

# New "aBc!2@#xYz" test block (grid in rows 7-10)
$ws.Range("C8").Value = "aBc!2@#xYz"
$ws.Range("R7").Value = "a"
$ws.Range("R9").Value = "c"
$ws.Range("R10").Value = "!"
$ws.Range("S10").Value = "@"
$ws.Range("T9").Value = "#"
$ws.Range("U8").Value = "x"
$ws.Range("V8").Value = "z"

# New "Hello World" test block (grid in rows 14-17)
$ws.Range("C9").Value = "Hello World"
$ws.Range("R15").Value = "e"
$ws.Range("R16").Value = "l"
$ws.Range("U14").Value = "W"
$ws.Range("U15").Value = "o"
$ws.Range("U16").Value = "r"
$ws.Range("V16").Value = "d"

# --- Remaining cells (values already interned above, any order) ---
$ws.Range("V7").Value = "Y"
$ws.Range("R8").Value = "B"
$ws.Range("R11").Value = 2
$ws.Range("R14").Value = "H"
$ws.Range("R17").Value = "l"
$ws.Range("U17").Value = "l"
$ws.Range("S16").Value = "o"

$ws.Range("H22").Value = "A"
$ws.Range("H23").Value = "A"
$ws.Range("K23").Value = "C"
$ws.Range("H24").Value = "A"
$ws.Range("I24").Value = "B"
$ws.Range("K24").Value = "C"
$ws.Range("L24").Value = "D"
$ws.Range("H25").Value = "B"
$ws.Range("K25").Value = "D"

# --- Selection matches the author's final cursor position ---
[void]$ws.Range("L16").Select()
